$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns at D:E (existing D:K shifts right to F:M)
$ws.Columns("D:E").Insert()

# Step 2: copy number formats / styles for the new D:E columns from the
#         (now-shifted) original D:E columns, which landed at F:G
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: populate the new D:E columns with the new quarter data
$newData = @{
    7 = @(43465, 43373)
    8 = @(6900, 5700)
    9 = @(1900, 2100)
    10 = @(5000, 3600)
    12 = @(1900, 1600)
    13 = @(0, 0)
    14 = @(-100, -5300)
    15 = @(200, 200)
    17 = @(7400, 2000)
    18 = @(-500, 3700)
    20 = @(600, -800)
    21 = @(500, 3400)
    22 = @(0, 0)
    23 = @(0, 2900)
    24 = @(0, 500)
    25 = @(0, 0)
    26 = @(0, 2400)
    27 = @(0, 2400)
    28 = @(0, 0)
    29 = @("NA", "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-600, 800)
    33 = @(0, 2400)
    34 = @(0, 0)
    35 = @(0, 2400)
    38 = @(43465, 43373)
    41 = @(8600, 8500)
    42 = @(0, 0)
    43 = @(7100, 4600)
    44 = @(0, 0)
    45 = @(2200, 1900)
    46 = @(17900, 15100)
    47 = @(0, 0)
    48 = @(500, 500)
    49 = @(11200, 11800)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(600, 1000)
    53 = @(0, 0)
    54 = @(30300, 28300)
    57 = @(2800, 1800)
    58 = @(200, 100)
    59 = @(14100, 13600)
    60 = @(17100, 15600)
    61 = @(3400, 3300)
    62 = @(2800, 2500)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(23300, 21400)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(-58900, -58900)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(7000, 7000)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(0, 2400)
    83 = @(400, 600)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(200, -200)
    91 = @(0, 0)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(100, 9600)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-100, -6100)
    101 = @(0, 0)
    102 = @(100, 3300)
}
foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item($r, 4).Value2 = $vals[0]
    $ws.Cells.Item($r, 5).Value2 = $vals[1]
}

# Step 4: a couple of cells among the shifted historical columns were also
#         corrected as part of this data refresh
$ws.Cells.Item(14, 6).Value2 = "NA"
$ws.Cells.Item(14, 7).Value2 = "NA"
$ws.Cells.Item(14, 8).Value2 = "NA"
$ws.Cells.Item(14, 9).Value2 = "NA"
$ws.Cells.Item(14, 10).Value2 = "NA"
$ws.Cells.Item(61, 10).Value2 = 6700

Write-Host "edit complete"
